$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 491.66666
$ws.Range("I12").Value = 900
$ws.Range("J12").Value = 287.5
$ws.Range("K12").Value = 900
$ws.Range("L12").Value = 287.5
$ws.Range("M12").Value = -730
$ws.Range("N12").Value = -627.5
$ws.Range("H18").Value = 627
$ws.Range("I18").Value = 627
$ws.Range("K18").Value = 627
$ws.Range("M18").Value = -343
$ws.Range("H40").Value = 1465.8334
$ws.Range("I40").Value = 1466.4706
$ws.Range("J40").Value = 1464.2858
$ws.Range("K40").Value = 1466.4706
$ws.Range("L40").Value = 1464.2858
$ws.Range("M40").Value = -1291.4706
$ws.Range("N40").Value = -1814.2858
$ws.Range("H62").Value = 1733.0555
$ws.Range("I62").Value = 1599
$ws.Range("J62").Value = 1784.6154
$ws.Range("K62").Value = 1599
$ws.Range("L62").Value = 1784.6154
$ws.Range("M62").Value = -975
$ws.Range("N62").Value = -3032.6154
$ws.Range("H64").Value = 34487180
$ws.Range("I64").Value = 55559950
$ws.Range("K64").Value = 55559950
$ws.Range("M64").Value = -55559702
$ws.Range("H65").Value = 1733.0555
$ws.Range("I65").Value = 1599
$ws.Range("J65").Value = 1784.6154
$ws.Range("K65").Value = 7995
$ws.Range("L65").Value = 8923.076999999999
$ws.Range("M65").Value = -4875
$ws.Range("N65").Value = -15163.077
$ws.Range("H67").Value = 34487180
$ws.Range("I67").Value = 55559950
$ws.Range("K67").Value = 55559950
$ws.Range("M67").Value = -55559092
$ws.Range("H69").Value = 4121.2383
$ws.Range("I69").Value = 5133.3335
$ws.Range("J69").Value = 3845.2122
$ws.Range("K69").Value = 15400.0005
$ws.Range("L69").Value = 11535.6366
$ws.Range("M69").Value = -14526.0005
$ws.Range("N69").Value = -13283.6366
$ws.Range("H70").Value = 3290
$ws.Range("I70").Value = 5340
$ws.Range("J70").Value = 1240
$ws.Range("K70").Value = 16020
$ws.Range("L70").Value = 3720
$ws.Range("M70").Value = -15750
$ws.Range("N70").Value = -4260
$ws.Range("H72").Value = 4121.2383
$ws.Range("I72").Value = 5133.3335
$ws.Range("J72").Value = 3845.2122
$ws.Range("K72").Value = 46200.0015
$ws.Range("L72").Value = 34606.9098
$ws.Range("M72").Value = -41832.0015
$ws.Range("N72").Value = -43342.9098
$ws.Range("H73").Value = 3290
$ws.Range("I73").Value = 5340
$ws.Range("J73").Value = 1240
$ws.Range("K73").Value = 16020
$ws.Range("L73").Value = 3720
$ws.Range("M73").Value = -15084
$ws.Range("N73").Value = -5592
$ws.Range("H74").Value = 3683.68
$ws.Range("I74").Value = 3443.3572
$ws.Range("J74").Value = 3989.5454
$ws.Range("K74").Value = 3443.3572
$ws.Range("L74").Value = 3989.5454
$ws.Range("M74").Value = -2507.3572
$ws.Range("N74").Value = -5861.5454
$ws.Range("H77").Value = 3683.68
$ws.Range("I77").Value = 3443.3572
$ws.Range("J77").Value = 3989.5454
$ws.Range("K77").Value = 17216.786
$ws.Range("L77").Value = 19947.727
$ws.Range("M77").Value = -12536.786
$ws.Range("N77").Value = -29307.727
$ws.Range("H80").Value = 570.8387
$ws.Range("I80").Value = 350.9091
$ws.Range("J80").Value = 691.8
$ws.Range("K80").Value = 1052.7273
$ws.Range("L80").Value = 2075.4
$ws.Range("M80").Value = -54.72730000000001
$ws.Range("N80").Value = -4071.4
$ws.Range("H83").Value = 570.8387
$ws.Range("I83").Value = 350.9091
$ws.Range("J83").Value = 691.8
$ws.Range("K83").Value = 3158.1819
$ws.Range("L83").Value = 6226.2
$ws.Range("M83").Value = 1833.8181
$ws.Range("N83").Value = -16210.2
$ws.Range("H100").Value = 1428.8387
$ws.Range("I100").Value = 1110.4762
$ws.Range("J100").Value = 2097.4
$ws.Range("K100").Value = 1110.4762
$ws.Range("L100").Value = 2097.4
$ws.Range("M100").Value = -569.4762000000001
$ws.Range("N100").Value = -3179.4
$ws.Range("H103").Value = 694
$ws.Range("I103").Value = 691
$ws.Range("J103").Value = 700
$ws.Range("K103").Value = 2073
$ws.Range("L103").Value = 2100
$ws.Range("M103").Value = -1487
$ws.Range("N103").Value = -3272
$ws.Range("H113").Value = 2788.3333
$ws.Range("I113").Value = 2096.5386
$ws.Range("J113").Value = 3912.5
$ws.Range("K113").Value = 2096.5386
$ws.Range("L113").Value = 3912.5
$ws.Range("M113").Value = 1157.4614
$ws.Range("N113").Value = -10420.5
$ws.Range("H116").Value = 3258.5
$ws.Range("I116").Value = 4051.25
$ws.Range("J116").Value = 2730
$ws.Range("K116").Value = 4051.25
$ws.Range("L116").Value = 2730
$ws.Range("M116").Value = -609.25
$ws.Range("N116").Value = -9614
$ws.Range("H138").Value = 2428.9119
$ws.Range("I138").Value = 1538.5264
$ws.Range("J138").Value = 3556.7334
$ws.Range("K138").Value = 4615.5792
$ws.Range("L138").Value = 10670.2002
$ws.Range("M138").Value = 524.4207999999999
$ws.Range("N138").Value = -20950.2002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 56237
$ws.Range("I2").Value = 111576.336
$ws.Range("J2").Value = 897.6667
$ws.Range("K2").Value = 111576.336
$ws.Range("L2").Value = 897.6667
$ws.Range("M2").Value = -111463.336
$ws.Range("N2").Value = -1123.6667
$ws.Range("H32").Value = 6254.0186
$ws.Range("I32").Value = 5092.614
$ws.Range("J32").Value = 11364.2
$ws.Range("K32").Value = 5092.614
$ws.Range("L32").Value = 11364.2
$ws.Range("M32").Value = -4805.614
$ws.Range("N32").Value = -11938.2
$ws.Range("H61").Value = 1023.84314
$ws.Range("I61").Value = 833.8837
$ws.Range("J61").Value = 2044.875
$ws.Range("K61").Value = 833.8837
$ws.Range("L61").Value = 2044.875
$ws.Range("M61").Value = -621.8837
$ws.Range("N61").Value = -2468.875
$ws.Range("H116").Value = 56237
$ws.Range("I116").Value = 111576.336
$ws.Range("J116").Value = 897.6667
$ws.Range("K116").Value = 111576.336
$ws.Range("L116").Value = 897.6667
$ws.Range("M116").Value = -109282.336
$ws.Range("N116").Value = -5485.6667
$ws.Range("H132").Value = 4420.0625
$ws.Range("I132").Value = 2879.5
$ws.Range("J132").Value = 11095.833
$ws.Range("K132").Value = 8638.5
$ws.Range("L132").Value = 33287.499
$ws.Range("M132").Value = -6108.5
$ws.Range("N132").Value = -38347.499
$ws.Range("H136").Value = 1023.84314
$ws.Range("I136").Value = 833.8837
$ws.Range("J136").Value = 2044.875
$ws.Range("K136").Value = 2501.6511
$ws.Range("L136").Value = 6134.625
$ws.Range("M136").Value = 48.34889999999996
$ws.Range("N136").Value = -11234.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 56237
$ws.Range("I3").Value = 111576.336
$ws.Range("J3").Value = 897.6667
$ws.Range("K3").Value = 111576.336
$ws.Range("L3").Value = 897.6667
$ws.Range("M3").Value = -111462.336
$ws.Range("N3").Value = -1125.6667
$ws.Range("H86").Value = 2594.9473
$ws.Range("I86").Value = 2470.6
$ws.Range("J86").Value = 2733.111
$ws.Range("K86").Value = 2470.6
$ws.Range("L86").Value = 2733.111
$ws.Range("M86").Value = -1347.6
$ws.Range("N86").Value = -4979.111
$ws.Range("H89").Value = 2594.9473
$ws.Range("I89").Value = 2470.6
$ws.Range("J89").Value = 2733.111
$ws.Range("K89").Value = 12353
$ws.Range("L89").Value = 13665.555
$ws.Range("M89").Value = -6737
$ws.Range("N89").Value = -24897.555
$ws.Range("H94").Value = 2683.843
$ws.Range("I94").Value = 545.2727
$ws.Range("J94").Value = 6604.5557
$ws.Range("K94").Value = 545.2727
$ws.Range("L94").Value = 6604.5557
$ws.Range("M94").Value = -94.27269999999999
$ws.Range("N94").Value = -7506.5557
$ws.Range("H105").Value = 1355.8334
$ws.Range("I105").Value = 1122.6364
$ws.Range("J105").Value = 1722.2858
$ws.Range("K105").Value = 1122.6364
$ws.Range("L105").Value = 1722.2858
$ws.Range("M105").Value = 624.3635999999999
$ws.Range("N105").Value = -5216.2858
$ws.Range("H134").Value = 754.6842
$ws.Range("I134").Value = 655.76086
$ws.Range("J134").Value = 1168.3636
$ws.Range("K134").Value = 1967.28258
$ws.Range("L134").Value = 3505.0908
$ws.Range("M134").Value = 567.7174199999999
$ws.Range("N134").Value = -8575.0908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1090.579
$ws.Range("I58").Value = 1020.6429
$ws.Range("J58").Value = 1286.4
$ws.Range("K58").Value = 1020.6429
$ws.Range("L58").Value = 1286.4
$ws.Range("M58").Value = -817.6429000000001
$ws.Range("N58").Value = -1692.4
$ws.Range("H105").Value = 1546.7778
$ws.Range("I105").Value = 1388.7142
$ws.Range("J105").Value = 2100
$ws.Range("K105").Value = 1388.7142
$ws.Range("L105").Value = 2100
$ws.Range("M105").Value = 358.2858000000001
$ws.Range("N105").Value = -5594
$ws.Range("H136").Value = 1090.579
$ws.Range("I136").Value = 1020.6429
$ws.Range("J136").Value = 1286.4
$ws.Range("K136").Value = 3061.9287
$ws.Range("L136").Value = 3859.2
$ws.Range("M136").Value = -511.9287000000004
$ws.Range("N136").Value = -8959.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 762.64514
$ws.Range("I132").Value = 765.5
$ws.Range("J132").Value = 761.2857
$ws.Range("K132").Value = 6889.5
$ws.Range("L132").Value = 6851.571300000001
$ws.Range("M132").Value = -4359.5
$ws.Range("N132").Value = -11911.5713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5953884
$ws.Range("I126").Value = 11905980
$ws.Range("J126").Value = 1787.75
$ws.Range("K126").Value = 35717940
$ws.Range("L126").Value = 5363.25
$ws.Range("M126").Value = -35715470
$ws.Range("N126").Value = -10303.25
$ws.Range("H132").Value = 3845.7358
$ws.Range("I132").Value = 4372.5
$ws.Range("J132").Value = 2224.923
$ws.Range("K132").Value = 13117.5
$ws.Range("L132").Value = 6674.768999999999
$ws.Range("M132").Value = -10587.5
$ws.Range("N132").Value = -11734.769
